$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the text of A2
$ws.Range("A2").Value = "Modification de la ligne 2 dans le dossier devoir *2"

# Remove row 3 entirely (deletes cell A3 and shifts cells up, but since nothing
# is below it, this simply clears out the row from the used range)
$ws.Range("A3").EntireRow.Delete() | Out-Null

# Restore the selection to A3 (now the first empty row below the data)
$ws.Range("A3").Select() | Out-Null
